# Updates the crypto price table (Sheet1) with the latest scraped prices /
# 1h volume percentages for the Wed Jun 5 23:33:45 UTC 2024 GitHub Actions run.
# Some rows also shift which coin occupies them (the source ranking reordered).
#
# Each entry below is: row number -> hashtable of column-letter -> new text value.
# D-column "price" values are plain decimal-look-alikes (e.g. "696.78"), so when
# written through Range.Value Excel's input parser would silently convert them
# to numbers; they must stay text (to match the original inlineStr cells), so
# we briefly force a Text number format, assign the value, then restore the
# cell's default style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    2  = @{ D = "70.879.55";  E = "  +0.62%  " }
    3  = @{ D = "3.855.16";   E = "  +1.20%  " }
    4  = @{ E = "  +0.02%  " }
    5  = @{ D = "696.78";     E = "  +2.26%  " }
    6  = @{ D = "173.17";     E = "  +1.08%  " }
    7  = @{ D = "3.853.23";   E = "  +1.19%  " }
    8  = @{ E = "  +0.02%  " }
    9  = @{ E = "  -0.01%  " }
    10 = @{ E = "  +1.04%  " }
    11 = @{ D = "7.16";       E = "  -1.82%  " }
    12 = @{ E = "  -0.13%  " }
    13 = @{ E = "  +4.84%  " }
    14 = @{ D = "36.40";      E = "  +1.15%  " }
    15 = @{ D = "4.504.27";   E = "  +1.19%  " }
    16 = @{ D = "3.860.05";   E = "  +1.31%  " }
    17 = @{ D = "71.038.46";  E = "  +0.72%  " }
    18 = @{ D = "17.71";      E = "  -0.01%  " }
    19 = @{ D = "7.24";       E = "  +0.77%  " }
    20 = @{ D = "0.115";      E = "  +0.13%  " }
    21 = @{ D = "11.18";      E = "  -0.78%  " }
    22 = @{ D = "491.99";     E = "  +3.11%  " }
    23 = @{ D = "0.724";      E = "  +1.35%  " }
    24 = @{ D = "85.01";      E = "  +1.92%  " }
    25 = @{ D = "0.0000145";  E = "  +1.76%  " }
    26 = @{ B = "RenderToken";               C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";              D = "10.72";     E = "  +3.71%  " }
    27 = @{ B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp";       D = "12.30";     E = "  +0.37%  " }
    28 = @{ E = "  +1.34%  " }
    29 = @{ B = "PancakeSwap";                C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake";                 D = "3.19";      E = "  +8.95%  " }
    30 = @{ B = "Dai";                        C = "https://coinranking.com/coin/MoTuySvg7+dai-dai";                          D = "1.00";      E = "  -0.05%  " }
    31 = @{ B = "NEARProtocol";               C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near";                D = "7.66";      E = "  +3.43%  " }
    32 = @{ B = "ImmutableX";                 C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";                   D = "2.28";      E = "  -0.60%  " }
    33 = @{ B = "EthereumClassic";            C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc";          D = "29.72";     E = "  +0.37%  " }
    34 = @{ B = "Kaspa";                      C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas";                        D = "0.178";     E = "  -1.23%  " }
    35 = @{ B = "Aptos";                      C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt";                        D = "9.29";      E = "  +1.67%  " }
    36 = @{ B = "RenzoRestakedETH";           C = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth";           D = "3.806.02";  E = "  +1.26%  " }
    37 = @{ B = "Binance-PegBSC-USD";         C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd";       D = "1.00";      E = "  +0.24%  " }
    38 = @{ B = "Hedera";                     C = "https://coinranking.com/coin/jad286TjB+hedera-hbar";                      D = "0.104";     E = "  +1.84%  " }
    39 = @{ B = "Stacks";                     C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx";                        D = "2.38";      E = "  +11.30%  " }
    40 = @{ B = "Filecoin";                   C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil";                     D = "6.07";      E = "  +1.94%  " }
    41 = @{ B = "dogwifhat";                  C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif";                    D = "3.39";      E = "  -0.15%  " }
    42 = @{ B = "Mantle";                     C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt";                       D = "1.03";      E = "  +7.29%  " }
    43 = @{ B = "FirstDigitalUSD";            C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd";            D = "1.00";      E = "  +0.15%  " }
    44 = @{ B = "USDe";                       C = "https://coinranking.com/coin/exbfr2U-0+usde-usde" }
    45 = @{ B = "Monero";                     C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr";                   D = "163.29";    E = "  +2.04%  " }
    46 = @{ B = "FLOKI";                      C = "https://coinranking.com/coin/fmHk13Rqw+floki-floki";                      D = "0.000308";  E = "  +4.38%  " }
    47 = @{ B = "OKB";                        C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb";                      D = "48.58";     E = "  +0.73%  " }
    48 = @{ B = "Arweave";                    C = "https://coinranking.com/coin/7XWg41D1+arweave-ar";                        D = "44.39";     E = "  -4.06%  " }
    49 = @{ B = "Bittensor";                  C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao";                    D = "419.61";    E = "  +5.24%  " }
    50 = @{ B = "TheGraph";                   C = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt";                     D = "0.303";     E = "  +1.14%  " }
    51 = @{ B = "Cosmos";                     C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom";                  D = "8.68";      E = "  +1.98%  " }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$rowNum"
        $newValue = $cols[$col]
        $cell = $ws.Range($cellRef)

        if ($col -eq "D") {
            # Force text storage for decimal-look-alike price strings so they
            # don't get auto-converted to numbers, then strip the temporary
            # Text format so the cell's style matches the untouched cells.
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newValue
        }
    }
}
